# Rework the FamilyMemberHistory -> FamilyHistory mapping sheet so the
# "FamilyHistory" column (B) lines up with Astrid's revised mapping.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Re-point the existing FamilyHistory.* values against the EHDS
#        rows they now map to (column A, rows 2-17, is untouched) ---
$ws.Range("B2").Value  = "FamilyHistory"
$ws.Range("B3").ClearContents()
$ws.Range("B4").ClearContents()
$ws.Range("B5").ClearContents()
$ws.Range("B6").ClearContents()
$ws.Range("B7").ClearContents()
$ws.Range("B8").Value  = "FamilyHistory.Date"
$ws.Range("B9").ClearContents()
$ws.Range("B10").ClearContents()
$ws.Range("B11").ClearContents()
$ws.Range("B12").ClearContents()
$ws.Range("B13").Value = "FamilyHistory.FamilyMember.BiologicalRelationship"
$ws.Range("B14").ClearContents()
$ws.Range("B15").Value = "FamilyHistory.FamilyMember.AgeAtDeath"
$ws.Range("B16").Value = "FamilyHistory.FamilyMember.Disorder.DisorderFamilyMember.Diagnosis"

# --- 2. The "unmapped" leftover FamilyHistory.* elements, listed below
#        the main table (rows 18-26) ---
$ws.Range("B18").Value = "FamilyHistory.FamilyMember"
$ws.Range("B19").Value = "FamilyHistory.FamilyMember.Comment"
$ws.Range("B20").Value = "FamilyHistory.FamilyMember.Disorder"
$ws.Range("B21").Value = "FamilyHistory.FamilyMember.Disorder.IsCauseOfDeath"
$ws.Range("B22").Value = "FamilyHistory.FamilyMember.DeathIndicator"
$ws.Range("B23").Value = "FamilyHistory.FamilyMember.Disorder.DisorderFamilyMember"
$ws.Range("B24").Value = "FamilyHistory.FamilyMember.Disorder.DisorderFamilyMember"
$ws.Range("B25").Value = "FamilyHistory.FamilyMember.AgeAtDeath"
$ws.Range("B26").Value = "FamilyHistory.FamilyMember.DeathIndicator"

# Rows 27-28 no longer exist in the reworked sheet.
$ws.Range("B27:B28").ClearContents()

# --- 3. Apply left/top alignment to the mapped cells (B2:B16 & B18:B22).
#        Format the source cell once, then copy/paste-special (formats
#        only) so only a single new style entry is created instead of
#        one per property assignment. ---
$src = $ws.Cells.Item(2, 2)
$src.HorizontalAlignment = -4131   # xlHAlignLeft
$src.VerticalAlignment = -4160     # xlVAlignTop
$src.Copy()
$ws.Range("B3:B16").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("B18:B22").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# --- 4. Widen column B to fit the longer mapping strings ---
$ws.Columns.Item(2).ColumnWidth = 46.998697916666664

# --- 5. Restore the cursor/selection as left by the author ---
$ws.Range("E25").Select()
